$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("I10").Value = 0.1067
$ws.Range("J10").Value = -0.0452
$ws.Range("K10").Value = 0.295
$ws.Range("L10").Value = 0.2573
$ws.Range("M10").Value = 0.3088
$ws.Range("N10").Value = 0.2796
$ws.Range("O10").Value = 0.1561
$ws.Range("P11").Value = -0.1221
$ws.Range("Q11").Value = -0.2406
$ws.Range("R11").Value = -0.226
$ws.Range("S11").Value = -0.2871
$ws.Range("T11").Value = -0.2221
$ws.Range("U11").Value = -0.1724
$ws.Range("V11").Value = -0.035
$ws.Range("W11").Value = -0.1145
$ws.Range("X11").Value = -0.4321
$ws.Range("I30").Value = -2.2874
$ws.Range("J30").Value = -0.5607
$ws.Range("K30").Value = 0.0638
$ws.Range("L30").Value = -0.1884
$ws.Range("M30").Value = 0.5287
$ws.Range("N30").Value = 0.2312
$ws.Range("O30").Value = -0.2545
$ws.Range("P31").Value = -0.7614
$ws.Range("Q31").Value = -0.4572
$ws.Range("R31").Value = -0.2208
$ws.Range("S31").Value = -0.8013
$ws.Range("T31").Value = -0.9297
$ws.Range("U31").Value = -0.6619
$ws.Range("V31").Value = -0.4547
$ws.Range("W31").Value = -0.1598
$ws.Range("X31").Value = -71.7788
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = 0
$ws.Range("N66").Value = 0
$ws.Range("O66").Value = 0
$ws.Range("P67").Value = -0.0018
$ws.Range("Q67").Value = -0.0024
$ws.Range("R67").Value = -0.0025
$ws.Range("S67").Value = -0.0032
$ws.Range("T67").Value = -0.0026
$ws.Range("U67").Value = -0.0017
$ws.Range("V67").Value = 0.0005
$ws.Range("W67").Value = -0.0004
$ws.Range("X67").Value = -0.0049
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("N86").Value = 0
$ws.Range("O86").Value = 0
$ws.Range("P87").Value = -0.0876
$ws.Range("Q87").Value = -0.0942
$ws.Range("R87").Value = -0.1018
$ws.Range("S87").Value = -0.1151
$ws.Range("T87").Value = -0.1102
$ws.Range("U87").Value = -0.0977
$ws.Range("V87").Value = -0.0843
$ws.Range("W87").Value = -0.0689
$ws.Range("X87").Value = -1.1317
